# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column E ("municipio-nombre") is re-curated from a measure to a
# dimension, matching how "provincia-nombre" (F) and "comarca-nombre" (K)
# are already modeled:
#   E2: iaest-measure:municipio-nombre -> sdmx-dimension:refArea
#   E3: medida                        -> dim
#   E4: xsd:int                       -> URI-Municipio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"
